{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block that was\n// appended to the bibliography section, along with the blank paragraph that\n// separates it from the last bibliography entry (FLEMMING, Diva M. ...).\nconst body = context.document.body;\n\nconst results = body.search(\"Ver no Jupiter Salvar em pdf Salvar em docx\", {\n  matchCase: true,\n});\ncontext.load(results, \"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not locate the 'Ver no Jupiter ...' paragraph\");\n}\n\n// The blank paragraph right before \"Ver no Jupiter ...\" and the copyright\n// paragraph right after it are removed together with it.\nconst jupiterParagraph = results.items[0].paragraphs.getFirst();\nconst blankParagraph = jupiterParagraph.getPrevious();\nconst copyrightParagraph = jupiterParagraph.getNext();\n\ncopyrightParagraph.delete();\njupiterParagraph.delete();\nblankParagraph.delete();\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block that was\n# appended to the bibliography section, along with the blank paragraph that\n# separates it from the last bibliography entry (FLEMMING, Diva M. ...).\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.Text = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not locate the 'Ver no Jupiter ...' paragraph\"\n}\n\n$jupiterParagraph = $searchRange.Paragraphs.Item(1)\n$jupiterIndex = $jupiterParagraph.Index\n\n# The blank paragraph right before \"Ver no Jupiter ...\" and the copyright\n# paragraph right after it are removed together with it.\n$blankParagraph = $d.Paragraphs.Item($jupiterIndex - 1)\n$copyrightParagraph = $d.Paragraphs.Item($jupiterIndex + 1)\n\n$deleteRange = $d.Range($blankParagraph.Range.Start, $copyrightParagraph.Range.End)\n$deleteRange.Delete()\n"}
